$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.731.56'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '2.251.38'
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'253.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').Value = "'0.635"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').Value = "'70.72"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.08%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = "'0.650"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.61%  '
$ws.Range('D10').Value = "'41.28"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.23%  '
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').Value = "'0.0962"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.96%  '
$ws.Range('D13').Value = "'7.36"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.54%  '
$ws.Range('D14').Value = "'0.104"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = '2.589.75'
$ws.Range('E15').Value = '  +3.34%  '
$ws.Range('D16').Value = "'0.890"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = "'14.85"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '2.263.35'
$ws.Range('E18').Value = '  +4.20%  '
$ws.Range('D19').Value = '42.681.33'
$ws.Range('E19').Value = '  +3.70%  '
$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').Value = "'6.26"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').Value = "'73.05"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.72%  '
$ws.Range('D23').Value = "'235.66"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.20%  '
$ws.Range('D25').Value = "'4.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('D26').Value = "'11.68"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.57%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  -4.09%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').Value = "'167.85"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').Value = "'21.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = "'6.12"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.47%  '
$ws.Range('E34').Value = '  +5.68%  '
$ws.Range('D35').Value = "'0.0791"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.54%  '
$ws.Range('D36').Value = "'0.125"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('D37').Value = "'28.05"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.18%  '
$ws.Range('D38').Value = "'4.70"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('D39').Value = "'4.20"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('E40').Value = '  +6.70%  '
$ws.Range('E41').Value = '  +3.67%  '
$ws.Range('D42').Value = "'12.61"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = "'5.84"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('D44').Value = "'64.21"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = "'0.204"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = "'4.97"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.31%  '
$ws.Range('D47').Value = "'8.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('E49').Value = '  +5.44%  '
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').Value = "'4.49"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.63%  '
